# Added Week 15 simulations
# Update Road ("R") row (row 3) target depth data on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# OFF sheet (row 3 = "R")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 418
$wsOff.Range("C3").Value = 284
$wsOff.Range("D3").Value = 89
$wsOff.Range("E3").Value = 38
$wsOff.Range("F3").Value = 6

# DEF sheet (row 3 = "R")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 448
$wsDef.Range("C3").Value = 303
$wsDef.Range("D3").Value = 106
$wsDef.Range("E3").Value = 52
$wsDef.Range("F3").Value = 6
$wsDef.Range("G3").Value = 6
